$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = "Aris"
$ws.Range("B32").Value = 45350
$ws.Range("B32").NumberFormat = "m/d/yy"
$ws.Range("C32").Formula = "=18+24/60"
$ws.Range("D32").Formula = "=20+21/60"

$ws.Range("F32").Select()
